$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 4
$ws.Range("B4").Value = "User clicks play on a video."
$ws.Range("C4").Value = "1. Navigate to a video`n2. Click the play button"
$ws.Range("D4").Value = "Video is played smoothly."
$ws.Range("E4").Value = "Pass"
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 38.25

# Row 5
$ws.Range("B5").Value = "User clicks the pause button on a video. User then resumes the video."
$ws.Range("C5").Value = "1. Navigate to a video`n2. Click the play button`n3. Click the pause button`n4. Click the play button again"
$ws.Range("D5").Value = "Video is paused. Resumes from the same timestamp it was paused."
$ws.Range("E5").Value = "Pass"
$ws.Range("C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 69

# Row 6
$ws.Range("B6").Value = "User drags the seek bar to a specific point in a video."
$ws.Range("C6").Value = "1. Navigate to a video`n2. Click and hold the drag bar and drag it to desired time stamp"
$ws.Range("D6").Value = "Video plays smoothly from the correct time stamp."
$ws.Range("E6").Value = "Pass"
$ws.Range("C6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 54.75

# Row 7
$ws.Range("B7").Value = "User clicks the fullscreen button on a video."
$ws.Range("C7").Value = "1. Navigate to a video`n2. Press the full screen button"
$ws.Range("D7").Value = "The video enters full screen mode. Exits full screen mode using ESC key."
$ws.Range("E7").Value = "Pass"
$ws.Range("C7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 38.25

# Row 8
$ws.Range("B8").Value = "User adjusts the volume."
$ws.Range("C8").Value = "1. Navigate to a video.`n2. Click and drag the audio to desired level"
$ws.Range("D8").Value = "Audio is changed accordingly."
$ws.Range("E8").Value = "Pass"
$ws.Range("C8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 50.25

# Row 9
$ws.Range("B9").Value = "User auto-plays a video."
$ws.Range("C9").Value = "1. Navigate to a video`n2. Finish the video`n3. Let the next video auto-play"
$ws.Range("D9").Value = "The next video starts playing automatically. "
$ws.Range("E9").Value = "Pass"
$ws.Range("C9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 54

# Row 10
$ws.Range("B10").Value = "User enables closed captioning."
$ws.Range("C10").Value = "1. Navigate to a video`n2. Click on the CC button"
$ws.Range("D10").Value = "Captions appear on the screen."
$ws.Range("E10").Value = "Pass"
$ws.Range("C10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 38.25

# Column width change. Target stored width is 26.5703125 characters; the COM
# ColumnWidth setter in this runtime rounds through an internal pixel grid, so
# the nearest value it can reproduce is 26.5. Empirically, an input of 25.67
# lands reliably in the bucket that stores back out as 26.5 (closest achievable
# approximation of the target 26.5703125).
$ws.Columns.Item(2).ColumnWidth = 25.67

# Selection change
$ws.Range("E11").Select()